$d = $word.ActiveDocument

$replacements = @(
    @("829×8=6632", "940×6=5640"),
    @("980×3=2940", "975×6=5850"),
    @("778×3=2334", "141×8=1128"),
    @("272×2=544", "770×5=3850"),
    @("910×5=4550", "267×6=1602"),
    @("756×9=6804", "682×6=4092"),
    @("445×2=890", "934×9=8406"),
    @("534×6=3204", "575×6=3450"),
    @("626×7=4382", "103×3=309"),
    @("880×7=6160", "788×6=4728"),
    @("518×2=1036", "369×9=3321"),
    @("245×9=2205", "822×6=4932"),
    @("239×7=1673", "727×6=4362"),
    @("944×7=6608", "463×5=2315"),
    @("738×6=4428", "347×7=2429"),
    @("402×4=1608", "769×9=6921"),
    @("991×3=2973", "266×9=2394"),
    @("925×5=4625", "336×9=3024"),
    @("168×6=1008", "988×3=2964"),
    @("869×3=2607", "355×4=1420"),
    @("412×4=1648", "383×3=1149"),
    @("608×9=5472", "320×4=1280"),
    @("941×4=3764", "720×8=5760"),
    @("798×8=6384", "169×3=507"),
    @("289×8=2312", "429×4=1716")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
